$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) cells must stay plain text so formats like "93.00"/"1.00" are preserved
# (Excel would otherwise coerce them to numbers and drop trailing zeros).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.527.95'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.481.41'
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.00'
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('E7').Value = '  -1.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.507'
$ws.Range('E9').Value = '  +2.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.75'
$ws.Range('E10').Value = '  -1.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0785'
$ws.Range('E11').Value = '  +1.03%  '
$ws.Range('E12').Value = '  +2.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.863.35'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.85'
$ws.Range('E14').Value = '  -1.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.21'
$ws.Range('E15').Value = '  +9.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.485.87'
$ws.Range('E16').Value = '  +1.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.768'
$ws.Range('E17').Value = '  -1.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.518.10'
$ws.Range('E18').Value = '  +1.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.41'
$ws.Range('E19').Value = '  +2.14%  '
$ws.Range('E20').Value = '  +2.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.60'
$ws.Range('E21').Value = '  +4.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.29'
$ws.Range('E22').Value = '  +1.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.61'
$ws.Range('E23').Value = '  +0.74%  '
$ws.Range('E24').Value = '  -1.59%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').Value = '  +4.55%  '
$ws.Range('E28').Value = '  +0.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.65'
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '158.22'
$ws.Range('E31').Value = '  +3.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.45'
$ws.Range('E32').Value = '  -0.25%  '
$ws.Range('E33').Value = '  +0.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0755'
$ws.Range('E34').Value = '  +1.99%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.36'
$ws.Range('E35').Value = '  +1.84%  '
$ws.Range('E36').Value = '  -8.45%  '
$ws.Range('E37').Value = '  +4.28%  '
$ws.Range('E38').Value = '  -3.67%  '
$ws.Range('E39').Value = '  -2.38%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('E41').Value = '  -1.41%  '
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.53'
$ws.Range('E43').Value = '  -3.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.979.54'
$ws.Range('E44').Value = '  +0.54%  '
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.96'
$ws.Range('E46').Value = '  -2.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.96'
$ws.Range('E47').Value = '  +3.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.720.91'
$ws.Range('E48').Value = '  +0.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '97.73'
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.12'
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.42'
$ws.Range('E51').Value = '  -1.72%  '
